$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vacation")
$ws.Select()
$ws.Range("I1").Value = "Type"
$ws.Range("I2").Value = "Annual Vacation"
